$wb = $excel.ActiveWorkbook

# --- 1) Insert a new worksheet '2022-Q1' right before the '总计' sheet ---
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Use an already-styled sheet as a formatting donor so the header row / index
# column pick up the same cell style (s=2: bold, centered, bordered) that every
# other quarter sheet in this workbook already uses.
$donor = $wb.Worksheets.Item("2021-Q4")
$donor.Range("B1:H1").Copy($q1.Range("B1:H1"))
$donor.Range("A2:A32").Copy($q1.Range("A2:A32"))

# --- 2) Header row ---
$q1.Cells.Item(1,2).Value = '基金代码'
$q1.Cells.Item(1,3).Value = '基金名称'
$q1.Cells.Item(1,4).Value = '基金规模'
$q1.Cells.Item(1,5).Value = '股票总仓位'
$q1.Cells.Item(1,6).Value = '仓位占比'
$q1.Cells.Item(1,7).Value = '持有市值(亿元)'
$q1.Cells.Item(1,8).Value = '仓位排名'

# --- 3) Fund rows (A already holds the right 0-based index from the donor copy) ---
# row 2: 163409
$q1.Cells.Item(2,2).Value = '''163409'
$q1.Cells.Item(2,3).Value = '兴全绿色投资混合(LOF)'
$q1.Cells.Item(2,4).Value = '''77.81'
$q1.Cells.Item(2,5).Value = '''89.66'
$q1.Cells.Item(2,6).Value = '''4.16'
$q1.Cells.Item(2,7).Value = '''3.2369'
$q1.Cells.Item(2,8).Value = 4
# row 3: 398051
$q1.Cells.Item(3,2).Value = '''398051'
$q1.Cells.Item(3,3).Value = '中海环保新能源混合'
$q1.Cells.Item(3,4).Value = '''22.97'
$q1.Cells.Item(3,5).Value = '''73.87'
$q1.Cells.Item(3,6).Value = '''4.52'
$q1.Cells.Item(3,7).Value = '''1.0382'
$q1.Cells.Item(3,8).Value = 8
# row 4: 005668
$q1.Cells.Item(4,2).Value = '''005668'
$q1.Cells.Item(4,3).Value = '融通新能源汽车主题精选灵活配置混合A'
$q1.Cells.Item(4,4).Value = '''22.12'
$q1.Cells.Item(4,5).Value = '''93.89'
$q1.Cells.Item(4,6).Value = '''3.26'
$q1.Cells.Item(4,7).Value = '''0.7211'
$q1.Cells.Item(4,8).Value = 10
# row 5: 519091
$q1.Cells.Item(5,2).Value = '''519091'
$q1.Cells.Item(5,3).Value = '新华泛资源优势混合'
$q1.Cells.Item(5,4).Value = '''13.39'
$q1.Cells.Item(5,5).Value = '''78.42'
$q1.Cells.Item(5,6).Value = '''2.71'
$q1.Cells.Item(5,7).Value = '''0.3629'
$q1.Cells.Item(5,8).Value = 9
# row 6: 000327
$q1.Cells.Item(6,2).Value = '''000327'
$q1.Cells.Item(6,3).Value = '南方潜力新蓝筹混合'
$q1.Cells.Item(6,4).Value = '''6.58'
$q1.Cells.Item(6,5).Value = '''93.61'
$q1.Cells.Item(6,6).Value = '''4.84'
$q1.Cells.Item(6,7).Value = '''0.3185'
$q1.Cells.Item(6,8).Value = 6
# row 7: 009885
$q1.Cells.Item(7,2).Value = '''009885'
$q1.Cells.Item(7,3).Value = '新华景气行业混合A'
$q1.Cells.Item(7,4).Value = '''9.84'
$q1.Cells.Item(7,5).Value = '''85.63'
$q1.Cells.Item(7,6).Value = '''2.83'
$q1.Cells.Item(7,7).Value = '''0.2785'
$q1.Cells.Item(7,8).Value = 8
# row 8: 519089
$q1.Cells.Item(8,2).Value = '''519089'
$q1.Cells.Item(8,3).Value = '新华优选成长混合'
$q1.Cells.Item(8,4).Value = '''8.58'
$q1.Cells.Item(8,5).Value = '''87.20'
$q1.Cells.Item(8,6).Value = '''3.05'
$q1.Cells.Item(8,7).Value = '''0.2617'
$q1.Cells.Item(8,8).Value = 9
# row 9: 519158
$q1.Cells.Item(9,2).Value = '''519158'
$q1.Cells.Item(9,3).Value = '新华趋势领航混合'
$q1.Cells.Item(9,4).Value = '''4.25'
$q1.Cells.Item(9,5).Value = '''89.82'
$q1.Cells.Item(9,6).Value = '''5.70'
$q1.Cells.Item(9,7).Value = '''0.2422'
$q1.Cells.Item(9,8).Value = 4
# row 10: 014150
$q1.Cells.Item(10,2).Value = '''014150'
$q1.Cells.Item(10,3).Value = '新华鑫益灵活配置混合A'
$q1.Cells.Item(10,4).Value = '''7.33'
$q1.Cells.Item(10,5).Value = '''84.57'
$q1.Cells.Item(10,6).Value = '''2.70'
$q1.Cells.Item(10,7).Value = '''0.1979'
$q1.Cells.Item(10,8).Value = 9
# row 11: 161605
$q1.Cells.Item(11,2).Value = '''161605'
$q1.Cells.Item(11,3).Value = '融通蓝筹成长混合'
$q1.Cells.Item(11,4).Value = '''4.82'
$q1.Cells.Item(11,5).Value = '''71.70'
$q1.Cells.Item(11,6).Value = '''3.37'
$q1.Cells.Item(11,7).Value = '''0.1624'
$q1.Cells.Item(11,8).Value = 7
# row 12: 000717
$q1.Cells.Item(12,2).Value = '''000717'
$q1.Cells.Item(12,3).Value = '融通转型三动力灵活配置混合A'
$q1.Cells.Item(12,4).Value = '''3.83'
$q1.Cells.Item(12,5).Value = '''94.89'
$q1.Cells.Item(12,6).Value = '''4.16'
$q1.Cells.Item(12,7).Value = '''0.1593'
$q1.Cells.Item(12,8).Value = 6
# row 13: 000584
$q1.Cells.Item(13,2).Value = '''000584'
$q1.Cells.Item(13,3).Value = '新华鑫益灵活配置混合'
$q1.Cells.Item(13,4).Value = '''5.89'
$q1.Cells.Item(13,5).Value = '''84.57'
$q1.Cells.Item(13,6).Value = '''2.70'
$q1.Cells.Item(13,7).Value = '''0.1590'
$q1.Cells.Item(13,8).Value = 9
# row 14: 519013
$q1.Cells.Item(14,2).Value = '''519013'
$q1.Cells.Item(14,3).Value = '海富通风格优势混合'
$q1.Cells.Item(14,4).Value = '''3.57'
$q1.Cells.Item(14,5).Value = '''89.46'
$q1.Cells.Item(14,6).Value = '''3.54'
$q1.Cells.Item(14,7).Value = '''0.1264'
$q1.Cells.Item(14,8).Value = 8
# row 15: 001280
$q1.Cells.Item(15,2).Value = '''001280'
$q1.Cells.Item(15,3).Value = '银华聚利灵活配置混合A'
$q1.Cells.Item(15,4).Value = '''4.58'
$q1.Cells.Item(15,5).Value = '''83.15'
$q1.Cells.Item(15,6).Value = '''2.72'
$q1.Cells.Item(15,7).Value = '''0.1246'
$q1.Cells.Item(15,8).Value = 10
# row 16: 012096
$q1.Cells.Item(16,2).Value = '''012096'
$q1.Cells.Item(16,3).Value = '鑫元鑫动力混合型证券投资基金A'
$q1.Cells.Item(16,4).Value = '''2.82'
$q1.Cells.Item(16,5).Value = '''88.33'
$q1.Cells.Item(16,6).Value = '''3.77'
$q1.Cells.Item(16,7).Value = '''0.1063'
$q1.Cells.Item(16,8).Value = 10
# row 17: 001105
$q1.Cells.Item(17,2).Value = '''001105'
$q1.Cells.Item(17,3).Value = '信达澳银转型创新股票'
$q1.Cells.Item(17,4).Value = '''2.54'
$q1.Cells.Item(17,5).Value = '''90.59'
$q1.Cells.Item(17,6).Value = '''3.67'
$q1.Cells.Item(17,7).Value = '''0.0932'
$q1.Cells.Item(17,8).Value = 8
# row 18: 011506
$q1.Cells.Item(18,2).Value = '''011506'
$q1.Cells.Item(18,3).Value = '建信高端装备股票型证券投资基金A'
$q1.Cells.Item(18,4).Value = '''2.59'
$q1.Cells.Item(18,5).Value = '''85.91'
$q1.Cells.Item(18,6).Value = '''3.55'
$q1.Cells.Item(18,7).Value = '''0.0919'
$q1.Cells.Item(18,8).Value = 10
# row 19: 009835
$q1.Cells.Item(19,2).Value = '''009835'
$q1.Cells.Item(19,3).Value = '融通新能源汽车主题精选灵活配置混合C'
$q1.Cells.Item(19,4).Value = '''2.28'
$q1.Cells.Item(19,5).Value = '''93.89'
$q1.Cells.Item(19,6).Value = '''3.26'
$q1.Cells.Item(19,7).Value = '''0.0743'
$q1.Cells.Item(19,8).Value = 10
# row 20: 011727
$q1.Cells.Item(20,2).Value = '''011727'
$q1.Cells.Item(20,3).Value = '工银瑞信聚瑞混合型证券投资基金A'
$q1.Cells.Item(20,4).Value = '''3.46'
$q1.Cells.Item(20,5).Value = '''29.56'
$q1.Cells.Item(20,6).Value = '''1.23'
$q1.Cells.Item(20,7).Value = '''0.0426'
$q1.Cells.Item(20,8).Value = 9
# row 21: 009886
$q1.Cells.Item(21,2).Value = '''009886'
$q1.Cells.Item(21,3).Value = '新华景气行业混合C'
$q1.Cells.Item(21,4).Value = '''1.28'
$q1.Cells.Item(21,5).Value = '''85.63'
$q1.Cells.Item(21,6).Value = '''2.83'
$q1.Cells.Item(21,7).Value = '''0.0362'
$q1.Cells.Item(21,8).Value = 8
# row 22: 009537
$q1.Cells.Item(22,2).Value = '''009537'
$q1.Cells.Item(22,3).Value = '太平行业优选股票A'
$q1.Cells.Item(22,4).Value = '''0.88'
$q1.Cells.Item(22,5).Value = '''90.50'
$q1.Cells.Item(22,6).Value = '''3.87'
$q1.Cells.Item(22,7).Value = '''0.0341'
$q1.Cells.Item(22,8).Value = 9
# row 23: 011507
$q1.Cells.Item(23,2).Value = '''011507'
$q1.Cells.Item(23,3).Value = '建信高端装备股票型证券投资基金C'
$q1.Cells.Item(23,4).Value = '''0.90'
$q1.Cells.Item(23,5).Value = '''85.91'
$q1.Cells.Item(23,6).Value = '''3.55'
$q1.Cells.Item(23,7).Value = '''0.0320'
$q1.Cells.Item(23,8).Value = 10
# row 24: 009828
$q1.Cells.Item(24,2).Value = '''009828'
$q1.Cells.Item(24,3).Value = '融通转型三动力灵活配置混合C'
$q1.Cells.Item(24,4).Value = '''0.59'
$q1.Cells.Item(24,5).Value = '''94.89'
$q1.Cells.Item(24,6).Value = '''4.16'
$q1.Cells.Item(24,7).Value = '''0.0245'
$q1.Cells.Item(24,8).Value = 6
# row 25: 014141
$q1.Cells.Item(25,2).Value = '''014141'
$q1.Cells.Item(25,3).Value = '大成新能源混合A'
$q1.Cells.Item(25,4).Value = '''0.52'
$q1.Cells.Item(25,5).Value = '''81.06'
$q1.Cells.Item(25,6).Value = '''3.93'
$q1.Cells.Item(25,7).Value = '''0.0204'
$q1.Cells.Item(25,8).Value = 9
# row 26: 004573
$q1.Cells.Item(26,2).Value = '''004573'
$q1.Cells.Item(26,3).Value = '新华鑫泰灵活配置混合'
$q1.Cells.Item(26,4).Value = '''0.73'
$q1.Cells.Item(26,5).Value = '''77.81'
$q1.Cells.Item(26,6).Value = '''2.56'
$q1.Cells.Item(26,7).Value = '''0.0187'
$q1.Cells.Item(26,8).Value = 8
# row 27: 002326
$q1.Cells.Item(27,2).Value = '''002326'
$q1.Cells.Item(27,3).Value = '银华聚利灵活配置混合C'
$q1.Cells.Item(27,4).Value = '''0.65'
$q1.Cells.Item(27,5).Value = '''83.15'
$q1.Cells.Item(27,6).Value = '''2.72'
$q1.Cells.Item(27,7).Value = '''0.0177'
$q1.Cells.Item(27,8).Value = 10
# row 28: 009538
$q1.Cells.Item(28,2).Value = '''009538'
$q1.Cells.Item(28,3).Value = '太平行业优选股票C'
$q1.Cells.Item(28,4).Value = '''0.20'
$q1.Cells.Item(28,5).Value = '''90.50'
$q1.Cells.Item(28,6).Value = '''3.87'
$q1.Cells.Item(28,7).Value = '''0.0077'
$q1.Cells.Item(28,8).Value = 9
# row 29: 002908
$q1.Cells.Item(29,2).Value = '''002908'
$q1.Cells.Item(29,3).Value = '富国睿利定期开放混合'
$q1.Cells.Item(29,4).Value = '''0.36'
$q1.Cells.Item(29,5).Value = '''28.52'
$q1.Cells.Item(29,6).Value = '''1.00'
$q1.Cells.Item(29,7).Value = '''0.0036'
$q1.Cells.Item(29,8).Value = 6
# row 30: 001866
$q1.Cells.Item(30,2).Value = '''001866'
$q1.Cells.Item(30,3).Value = '北信瑞丰新成长灵活配置混合'
$q1.Cells.Item(30,4).Value = '''0.07'
$q1.Cells.Item(30,5).Value = '''94.21'
$q1.Cells.Item(30,6).Value = '''4.69'
$q1.Cells.Item(30,7).Value = '''0.0033'
$q1.Cells.Item(30,8).Value = 7
# row 31: 011728
$q1.Cells.Item(31,2).Value = '''011728'
$q1.Cells.Item(31,3).Value = '工银瑞信聚瑞混合型证券投资基金C'
$q1.Cells.Item(31,4).Value = '''0.17'
$q1.Cells.Item(31,5).Value = '''29.56'
$q1.Cells.Item(31,6).Value = '''1.23'
$q1.Cells.Item(31,7).Value = '''0.0021'
$q1.Cells.Item(31,8).Value = 9
# row 32: 014142
$q1.Cells.Item(32,2).Value = '''014142'
$q1.Cells.Item(32,3).Value = '大成新能源混合C'
$q1.Cells.Item(32,4).Value = '''0.02'
$q1.Cells.Item(32,5).Value = '''81.06'
$q1.Cells.Item(32,6).Value = '''3.93'
$q1.Cells.Item(32,7).Value = '''0.0008'
$q1.Cells.Item(32,8).Value = 9

# --- 4) Update the '总计' summary sheet: insert a new top data row for 2022-Q1 ---
$tot = $wb.Worksheets.Item("总计")
$tot.Rows.Item(2).Insert()
$tot.Range("B2:D2").ClearFormats()

# Give the new A2 the same index-column style (s=2) already used by A3:A5
$tot.Cells.Item(3,1).Copy($tot.Cells.Item(2,1))

$tot.Cells.Item(2,1).Value = 0
$tot.Cells.Item(2,2).Value = "2022-Q1"
$tot.Cells.Item(2,3).Value = 31
$tot.Cells.Item(2,4).Value = 8

# Renumber the pre-existing rows' index column (they shifted down by one row)
$tot.Cells.Item(3,1).Value = 1
$tot.Cells.Item(4,1).Value = 2
$tot.Cells.Item(5,1).Value = 3

